# Add a "Results URL" column (new column E) to the meta-analysis dataset
# table, inserting a dropbox results-pickle link + hyperlink for each
# dataset row, and growing the row heights to match the now-wrapped
# extra column of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column E ("Description" et al
# shift one column to the right, from E..I to F..J).
$ws.Columns("E:E").Insert()

# Header for the new column.
$ws.Range("E1").Value = "Results URL"

# Per-row Results URL values (dropbox pickle links with the analysis
# results for each dataset), kept in variables so the same literal text
# is used both for the cell value and for the hyperlink target.
$uImsdb       = "https://www.dropbox.com/scl/fi/3gq5ieq7l25719if3my1f/imsdb_results.pkl?rlkey=01fjsk43sb8g05ccioysj0a7i&dl=1"
$uMovies      = "https://www.dropbox.com/scl/fi/arxkyhub2fi6qh5t79pfi/movies_results.pkl?rlkey=wmpf6aufzd2q86yju990a9keo&dl=1"
$uSwitchboard = "https://www.dropbox.com/scl/fi/1o7wqdlc1oo26y6ldpv8i/switchboard_results.pkl?rlkey=fetrrcp0vbsrmwydh39ikb918&dl=1"
$uScotus      = "https://www.dropbox.com/scl/fi/zxkvlrg4lfxcv7cjythp5/scotus_results.pkl?rlkey=krllpoa2jxvjlxrjxz6v9z2p5&dl=1"
$uTennis      = "https://www.dropbox.com/scl/fi/d3g83mtz4mqhbpmxfco5t/tennis_results.pkl?rlkey=ti9lsz49zyv8ru77cn2240qk5&dl=1"
$uPfg         = "https://www.dropbox.com/scl/fi/zmumd8uno58cqzoptr08m/pfg_results.pkl?rlkey=b4n8b7nh92rwgo7s91hgj7087&dl=1"
$uIq2         = "https://www.dropbox.com/scl/fi/3d4eha6r6xop7h0u1shgg/iq2_results.pkl?rlkey=qhaltntbg03len7bqqrwcpgjw&dl=1"
$uGap         = "https://www.dropbox.com/scl/fi/prk03sodn4pg8954cx9pa/gap_results.pkl?rlkey=wj2mngehnrm52thoetrkw22u1&dl=1"
$uChair       = "https://www.dropbox.com/scl/fi/to0642t939pvrtz1tka9y/chair_results.pkl?rlkey=sqz65t6sap29fkedwd7vwk3w6&dl=1"
$uFriends     = "https://www.dropbox.com/scl/fi/mkxc114g90rifsmzm881f/friends_results.pkl?rlkey=53qf44bwl2668h4irz14bf1ig&dl=1"
$uGutenberg   = "https://www.dropbox.com/s/jz15wcsceacaqva/gutenberg_results.pkl?dl=1"
$uReddit      = "https://www.dropbox.com/scl/fi/p999uknmzhx6f49ps5j0l/reddit_results.pkl?rlkey=lcmyva4t65dapqvw4y8fq8qrg&dl=1"

$ws.Range("E2").Value  = $uImsdb
$ws.Range("E3").Value  = $uMovies
$ws.Range("E4").Value  = $uSwitchboard
$ws.Range("E5").Value  = $uScotus
$ws.Range("E6").Value  = $uTennis
$ws.Range("E7").Value  = $uPfg
$ws.Range("E8").Value  = $uIq2
$ws.Range("E9").Value  = $uGap
$ws.Range("E10").Value = $uChair
$ws.Range("E11").Value = $uFriends
$ws.Range("E12").Value = $uGutenberg
$ws.Range("E13").Value = $uReddit

# Turn each of those URLs into a real hyperlink (target == cell text).
$ws.Hyperlinks.Add($ws.Range("E2"),  $uImsdb)
$ws.Hyperlinks.Add($ws.Range("E3"),  $uMovies)
$ws.Hyperlinks.Add($ws.Range("E4"),  $uSwitchboard)
$ws.Hyperlinks.Add($ws.Range("E5"),  $uScotus)
$ws.Hyperlinks.Add($ws.Range("E6"),  $uTennis)
$ws.Hyperlinks.Add($ws.Range("E7"),  $uPfg)
$ws.Hyperlinks.Add($ws.Range("E8"),  $uIq2)
$ws.Hyperlinks.Add($ws.Range("E9"),  $uGap)
$ws.Hyperlinks.Add($ws.Range("E10"), $uChair)
$ws.Hyperlinks.Add($ws.Range("E11"), $uFriends)
$ws.Hyperlinks.Add($ws.Range("E12"), $uGutenberg)
$ws.Hyperlinks.Add($ws.Range("E13"), $uReddit)

# Hyperlinks.Add re-styles the cells with the built-in "Hyperlink" look;
# restore the same (non-hyperlink-blue) formatting the rest of the table
# uses by re-pasting the formats from column D.
$ws.Range("D1:D13").Copy()
$ws.Range("E1:E13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column E should be as wide as column D (18 chars).
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth

# Growing the description text into a 6th visible column means several
# rows need to be a bit taller to keep everything wrapped & visible.
$ws.Rows("2:2").RowHeight = 102
$ws.Rows("3:3").RowHeight = 102
$ws.Rows("4:4").RowHeight = 119
$ws.Rows("5:5").RowHeight = 102
$ws.Rows("6:6").RowHeight = 102
$ws.Rows("7:7").RowHeight = 119
$ws.Rows("8:8").RowHeight = 102
$ws.Rows("9:9").RowHeight = 153
$ws.Rows("10:10").RowHeight = 102
$ws.Rows("11:11").RowHeight = 119
$ws.Rows("12:12").RowHeight = 68
$ws.Rows("13:13").RowHeight = 102

# Keep the sheet's remembered sort-range in sync with the inserted column
# (was A2:F13 / F2:F13, now one column wider: A2:G13 / G2:G13).
$ws.Range("A2:G13").Sort($ws.Range("G2:G13"), 1)

# Clear the old single-cell selection the workbook had remembered.
$ws.Range("A1").Select()
